$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before old column I (RATE)
$ws.Columns("I:J").Insert()

# Update header texts for the shifted-in columns
$ws.Range("G1").Value = "GSTR 3B FILED"
$ws.Range("H1").Value = "GSTR 1 FILED"
$ws.Range("I1").Value = "GSTR 1 FILLING DT"
$ws.Range("J1").Value = "GSTR 1 PERIOD"

# Apply unified header formatting across the whole header row
$hdr = $ws.Range("A1:Q1")
$hdr.NumberFormat = "0"
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 8
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.WrapText = $true
$ws.Rows(1).RowHeight = 33.75

# Row 2 updates
$ws.Range("D2").Clear()
$ws.Range("I2").Clear()
$ws.Range("H2").Value = "Y/N"
$ws.Range("K2").Clear()
$ws.Range("J2").Value = "Apr-19"
$ws.Range("J2").Font.Name = "Arial"
$ws.Range("J2").Font.Size = 10
$ws.Range("J2").HorizontalAlignment = -4152
$ws.Range("J2").NumberFormat = "@"

# Remove trailing empty rows (24-26)
$ws.Rows("24:26").Delete()
